# working on import issues
# Fix country values in Personnel sheet (column G) to use short codes
# instead of full country names, and update sheet selections to match
# the state the workbook was left in.

$wb = $excel.ActiveWorkbook

# --- Personnel sheet: normalize country codes in column G ---
$wsPersonnel = $wb.Worksheets.Item("Personnel")
$wsPersonnel.Range("G2:G16").Value = "TR"
$wsPersonnel.Range("G17:G19").Value = "GR"
$wsPersonnel.Range("G20").Value = "US"
$wsPersonnel.Range("G21").Value = "UK"
$wsPersonnel.Range("G22").Value = "DE"
$wsPersonnel.Range("G23").Value = "FR"
$wsPersonnel.Range("G24").Value = "IT"
$wsPersonnel.Range("G25").Value = "SP"
$wsPersonnel.Range("G26").Value = "PL"

# --- Countries sheet: leave selection on C4 ---
$wsCountries = $wb.Worksheets.Item("Countries")
$wsCountries.Activate()
$wsCountries.Range("C4").Select()

# --- Personnel sheet: finish with this sheet active and G2:G26 selected ---
$wsPersonnel.Activate()
$wsPersonnel.Range("G2:G26").Select()
